$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at row 577 (shifts old rows 577-595 down to 578-596) ---
$ws.Rows.Item(577).Insert()

# Populate the newly inserted row 577 with the new record
$ws.Cells.Item(577, 1).Value = 7
$ws.Cells.Item(577, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(577, 3).Value = "Ñuble"
$ws.Cells.Item(577, 4).Value = 45075
$ws.Cells.Item(577, 5).Value = 16
$ws.Cells.Item(577, 6).Value = 100114001
$ws.Cells.Item(577, 7).Value = "Papa"
$ws.Cells.Item(577, 8).Value = "Asterix"
$ws.Cells.Item(577, 9).Value = "1a (guarda lavada)"
$ws.Cells.Item(577, 10).Value = 100
$ws.Cells.Item(577, 11).Value = 12000
$ws.Cells.Item(577, 12).Value = 12000
$ws.Cells.Item(577, 13).Value = 12000
$ws.Cells.Item(577, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(577, 15).Value = "Región de Los Lagos"
$ws.Cells.Item(577, 16).Value = 480
$ws.Cells.Item(577, 17).Value = 25
$ws.Cells.Item(577, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Cells.Item(577, 4).NumberFormat = $ws.Cells.Item(578, 4).NumberFormat

# --- Append a brand-new row 597 at the bottom of the table ---
$ws.Cells.Item(597, 1).Value = 7
$ws.Cells.Item(597, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(597, 3).Value = "Ñuble"
$ws.Cells.Item(597, 4).Value = 44910
$ws.Cells.Item(597, 5).Value = 16
$ws.Cells.Item(597, 6).Value = 100114001
$ws.Cells.Item(597, 7).Value = "Papa"
$ws.Cells.Item(597, 8).Value = "Rosara"
$ws.Cells.Item(597, 9).Value = "1a nueva(o)"
$ws.Cells.Item(597, 10).Value = 160
$ws.Cells.Item(597, 11).Value = 11000
$ws.Cells.Item(597, 12).Value = 12000
$ws.Cells.Item(597, 13).Value = 11500
$ws.Cells.Item(597, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(597, 15).Value = "Región del Maule"
$ws.Cells.Item(597, 16).Value = 460
$ws.Cells.Item(597, 17).Value = 25
$ws.Cells.Item(597, 18).Value = "Hortaliza"

# Match the date number format used by the other rows in column D
$ws.Cells.Item(597, 4).NumberFormat = $ws.Cells.Item(596, 4).NumberFormat
